$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their text formatting so numeric-looking
# strings (e.g. "1.001", "1.400") are not coerced into numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.238.69"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.855.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7018"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "237.88"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07919"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3027"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.42"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08154"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.887.24"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.209"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7065"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.42"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.328.70"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.798"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007845"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.74"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.133.54"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.582"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.35"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.896"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1421"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.03"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.909"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.481"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.281"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.023"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05166"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.179"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.70%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7101"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.689"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.143.53"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +3.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9189"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.959"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +1.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4232"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.08%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "103.03"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5295"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.037.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.750"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.175"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.61%  "
